# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") values are recalculated/regenerated; write the new values in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 3
    4  = 0
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 6
    13 = 1
    14 = 5
    15 = 1
    16 = 3
    17 = 0
    18 = 2
    19 = 1
    20 = 4
    21 = 1
    22 = 4
    23 = 2
    24 = 2
    25 = 4
    26 = 4
    27 = 4
    28 = 4
    29 = 6
    30 = 2
    31 = 3
    32 = 1
    33 = 7
    34 = 0
    35 = 4
    36 = 5
    37 = 6
    38 = 4
    39 = 3
    40 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
